# Update "想去人数" (F column) figures across the four sheets to reflect
# the latest generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 29
$ws1.Range("F3").Value  = 1201
$ws1.Range("F5").Value  = 69
$ws1.Range("F7").Value  = 942
$ws1.Range("F8").Value  = 349
$ws1.Range("F9").Value  = 601
$ws1.Range("F11").Value = 1414
$ws1.Range("F12").Value = 130
$ws1.Range("F15").Value = 378
$ws1.Range("F17").Value = 1349
$ws1.Range("F18").Value = 777
$ws1.Range("F19").Value = 228
$ws1.Range("F20").Value = 1344
$ws1.Range("F23").Value = 1111
$ws1.Range("F25").Value = 3420
$ws1.Range("F26").Value = 669
$ws1.Range("F28").Value = 1518

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 46
$ws2.Range("F8").Value  = 18
$ws2.Range("F12").Value = 70

# --- Sheet: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 789

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 29
$ws4.Range("F3").Value  = 789
$ws4.Range("F4").Value  = 1201
$ws4.Range("F7").Value  = 69
$ws4.Range("F8").Value  = 46
$ws4.Range("F9").Value  = 46
$ws4.Range("F15").Value = 18
$ws4.Range("F17").Value = 942
$ws4.Range("F18").Value = 349
$ws4.Range("F19").Value = 601
$ws4.Range("F21").Value = 1414
$ws4.Range("F22").Value = 130
$ws4.Range("F25").Value = 378
$ws4.Range("F27").Value = 1349
$ws4.Range("F28").Value = 777
$ws4.Range("F29").Value = 228
$ws4.Range("F30").Value = 1344
$ws4.Range("F35").Value = 1111
$ws4.Range("F37").Value = 3420
$ws4.Range("F38").Value = 669
$ws4.Range("F40").Value = 1518
$ws4.Range("F41").Value = 70
